# Refitting NCDEs to individual patients (for manuscript figure)
# Add a new "Label" column (H) to the worksheet, mirroring the style of the
# other header cells (B1:G1) and filling in the 0/1 label values for each
# of the two 10-row blocks (rows 2-11 and rows 12-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, styled like the rest of row 1 (copy formatting from G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Label"

# First block (rows 2-11)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1

# Second block (rows 12-21)
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
